$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (tab) to reflect new "through" date
$ws.Name = "Through 2022-10-13"

# Update the header label in I1 (shared string text)
$ws.Range("I1").Value = "2022 (through 10-13)"

# Update the updated data values
$ws.Range("I11").Value = 44
$ws.Range("I14").Value = 1322
